# Weekly update: insert a new price record for Ajo (Macroferia Regional de
# Talca) as the first data row of this sub-range, pushing the existing
# records (rows 468-551) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(468).Insert()

$ws.Range("A468").Value = 5
$ws.Range("B468").Value = "Macroferia Regional de Talca"
$ws.Range("C468").Value = "Maule"
$ws.Range("D468").Value = 45209
$ws.Range("E468").Value = 7
$ws.Range("F468").Value = 100112003
$ws.Range("G468").Value = "Ajo"
$ws.Range("H468").Value = "Chino"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 200
$ws.Range("K468").Value = 20000
$ws.Range("L468").Value = 20000
$ws.Range("M468").Value = 20000
$ws.Range("N468").Value = "$/malla 10 kilos"
$ws.Range("O468").Value = "China"
$ws.Range("P468").Value = 2000
$ws.Range("Q468").Value = 10
$ws.Range("R468").Value = "Hortaliza"
